$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new user records (rows 6-8) to the Users_db sheet, matching the
# existing ID / Password / Type column layout (columns A / B / C).
$ws.Cells.Item(6, 1).Value = 258741359
$ws.Cells.Item(6, 2).Value = 15478
$ws.Cells.Item(6, 3).Value = 3

$ws.Cells.Item(7, 1).Value = 123852146
$ws.Cells.Item(7, 2).Value = 123
$ws.Cells.Item(7, 3).Value = 2

$ws.Cells.Item(8, 1).Value = 125478524
$ws.Cells.Item(8, 2).Value = 123
$ws.Cells.Item(8, 3).Value = 2
